# This script applies a permutation of the data in columns D, J, K, L, M, O, P
# across rows 2-21 of the active worksheet. Row r's new values come from the
# original values found in row mapping[r].

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (1-indexed, matches spreadsheet rows)
$mapping = @{
    2  = 16
    3  = 6
    4  = 2
    5  = 4
    6  = 5
    7  = 14
    8  = 15
    9  = 18
    10 = 7
    11 = 19
    12 = 17
    13 = 10
    14 = 12
    15 = 11
    16 = 3
    17 = 20
    18 = 9
    19 = 8
    20 = 13
    21 = 21
}

$cols = @("D", "J", "K", "L", "M", "O", "P")

# Snapshot the original values for the columns involved, for every row 2..21
$original = @{}
for ($r = 2; $r -le 21; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $original[$r] = $rowVals
}

# Apply the permutation using the snapshot so we never read already-overwritten data
for ($r = 2; $r -le 21; $r++) {
    $srcRow = $mapping[$r]
    $srcVals = $original[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $srcVals[$c]
    }
}
